$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Rule" label in B11 (row for rule R40) is renamed from "R40" to "1".
# A scratch cell is used to preserve B11's original cell formatting (style),
# because writing a numeric-looking string ("1") with a leading apostrophe
# (the only reliable way to force Excel to keep it as text instead of
# auto-converting it to the number 1) also stamps a "quote prefix" format
# onto the cell; copying the original format back over it afterwards keeps
# the cell's look identical to before the edit.
$target = $ws.Range("B11")
$scratch = $ws.Range("Z100")

$target.Copy($scratch)
$target.Value = "'1"
$scratch.Copy()
$target.PasteSpecial(-4122)
$scratch.Clear()
